$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Insert a new column C ("indicators_matched") into the "indicators" sheet.
#    First, duplicate column B's formatting into column C for rows 1-26,
#    then duplicate column B's values into column C (rows 1-26). This
#    reproduces both the per-row style (s="1" normal rows, s="2" highlighted
#    "Key Biodiversity Area" rows) and a starting value equal to column B.
# ---------------------------------------------------------------------------
$ws1.Range("B1:B26").Copy()
$ws1.Range("C1:C26").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("B1:B26").Copy()
$ws1.Range("C1:C26").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# Column C header
$ws1.Range("C1").Value2 = "indicators_matched"

# New descriptive text for the three "red list index" rows (differs from
# column B for these rows only; all other rows keep column B's value as-is)
$ws1.Range("C9").Value2  = "red list index of species survival for cms-listed bird species"
$ws1.Range("C10").Value2 = "red list index of species survival for migratory bird species"
$ws1.Range("C11").Value2 = "red list index of species survival for cms-listed bird and mammal species"

# Set the width of the new column C
$ws1.Columns.Item(3).ColumnWidth = 36

# ---------------------------------------------------------------------------
# 2) Rename column B header from "indicators" to "indicator"
# ---------------------------------------------------------------------------
$ws1.Range("B1").Value2 = "indicator"
